$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section 1 "Planificación" -------------------------------------------
# Row 2/3: drop the leading indentation spaces from the existing items.
$ws.Cells.Item(2, 1).Value = "1.1 Recolección de requisitos"
$ws.Cells.Item(3, 1).Value = "1.2 Análisis de viabilidad"

# Two brand-new sub-tasks get added right after them.
$ws.Rows.Item(4).Insert()
$ws.Cells.Item(4, 1).Value = "1.3 Asignación de roles"

$ws.Rows.Item(5).Insert()
$ws.Cells.Item(5, 1).Value = "1.4 Asignación de tareas"

# --- Section 2 "Diseño" ----------------------------------------------------
# (row 7 "2. Diseño" header, row 8 "   2.1 Diseño UI/UX" stay as-is)
# Two new sub-items of 2.1 are inserted into column B right after it.
$ws.Rows.Item(9).Insert()
$ws.Cells.Item(9, 2).Value = "2.1.1 Creación de mock-ups"

$ws.Rows.Item(10).Insert()
$ws.Cells.Item(10, 2).Value = "2.1.2 Creación y definición de estilos css y tailwind"

# (row 11 "   2.2 Arquitectura de software" stays as-is)

# The blank separator row that used to sit between "2. Diseño" and
# "3. Desarrollo" is removed - the sections now run back to back.
$ws.Rows.Item(12).Delete()

# --- Section 3 "Desarrollo" -------------------------------------------------
# (row 12 "3. Desarrollo" header stays as-is)
# Row 13 "3.1 Backend" loses its parenthetical description.
$ws.Cells.Item(13, 1).Value = "   3.1 Backend"

# Three new sub-items of 3.1 are inserted into column B right after it.
$ws.Rows.Item(14).Insert()
$ws.Cells.Item(14, 2).Value = "3.1.1 Creación de la base de datos"

$ws.Rows.Item(15).Insert()
$ws.Cells.Item(15, 2).Value = "3.1.2 Creación de las APIS"

$ws.Rows.Item(16).Insert()
$ws.Cells.Item(16, 2).Value = "3.1.3 Integración apis con el backend"

# Row 17 "3.2 Frontend" loses its parenthetical description (keeps the
# trailing space, no longer has a description after it).
$ws.Cells.Item(17, 1).Value = "   3.2 Frontend "

# --- Sections 4 "Pruebas" and 5 "Despliegue" --------------------------------
# Everything below is unchanged content-wise; it simply rode down with the
# row inserts/deletes above.

# --- Un-merge the section header bands (A:B) except the first one ----------
# The section headers for "2. Diseño", "3. Desarrollo", "4. Pruebas" and
# "5. Despliegue" are no longer merged across A:B (only "1. Planificación"
# keeps its merge). Unmerging leaves the fill/bold formatting on both cells.
$ws.Range("A7:B7").UnMerge() | Out-Null
$ws.Range("A12:B12").UnMerge() | Out-Null
$ws.Range("A19:B19").UnMerge() | Out-Null
$ws.Range("A24:B24").UnMerge() | Out-Null

# --- View state --------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 70
$ws.Range("H20").Select()
